$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(91).Insert()

$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 44729
$ws.Range("D91").NumberFormat = $ws.Range("D92").NumberFormat
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = 100114007
$ws.Range("G91").Value = "Jengibre"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 45
$ws.Range("K91").Value = 31000
$ws.Range("L91").Value = 31000
$ws.Range("M91").Value = 31000
$ws.Range("N91").Value = "`$/caja 13 kilos"
$ws.Range("O91").Value = "Perú"
$ws.Range("P91").Value = 2385
$ws.Range("Q91").Value = 13
$ws.Range("R91").Value = "Hortaliza"
